# Replace the single extracted-result data row (row 2) with a new record,
# and add the MiddleName column (I2) that the new record includes.
#
# NumberFormat is forced to "@" (text) immediately before each write and
# cleared again immediately after, so that date-like ("2024-02-06") and
# purely-numeric ("40402") values are stored as literal text instead of
# being auto-converted by Excel into date serials / numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

Set-TextValue "B2" "2024-02-06"
Set-TextValue "C2" "40402"
Set-TextValue "D2" "17500040"
Set-TextValue "E2" "16549"
Set-TextValue "F2" "Sagis DX"
Set-TextValue "G2" "TRESCH JR"
Set-TextValue "H2" "ROBERT"
Set-TextValue "I2" "E"
Set-TextValue "J2" "1952-04-24"
Set-TextValue "K2" "Male"
Set-TextValue "L2" "7104 DOSWELL LN"
Set-TextValue "M2" "TX"
Set-TextValue "N2" "AUSTIN"
Set-TextValue "O2" "787392042"
Set-TextValue "P2" "5124238633"
Set-TextValue "Q2" "table"
Set-TextValue "R2" "Leigha Sharp, MD,"
Set-TextValue "S2" "SWDV- Austin"
Set-TextValue "W2" "Medicare of Texas"
